# Update "想去人数" (F column) figures (and one status cell G) on the
# "展览" and "全部类型" worksheets, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId=1, first sheet) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 307
$ws1.Range("F7").Value = 2109
$ws1.Range("F10").Value = 4659
$ws1.Range("F16").Value = 150
$ws1.Range("F20").Value = 3541
$ws1.Range("F22").Value = 569
$ws1.Range("F32").Value = 754
$ws1.Range("F33").Value = 2193
$ws1.Range("F34").Value = 405
$ws1.Range("G34").Value = 128

# --- Sheet "全部类型" (sheetId=4, fourth sheet) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 307
$ws4.Range("F7").Value = 2109
$ws4.Range("F10").Value = 4659
$ws4.Range("F16").Value = 150
$ws4.Range("F20").Value = 3541
$ws4.Range("F22").Value = 569
$ws4.Range("F33").Value = 754
$ws4.Range("F34").Value = 2193
$ws4.Range("F35").Value = 405
$ws4.Range("G35").Value = 128
